$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "the number of first layer change from 8 to 16, Dropout from 0.2 to 0.15, epochs fromm 100 to 70"
$ws.Range("B7").Value = "new_arch_v2"
$ws.Range("G7").Value = 0.84210526332300395
$ws.Range("H7").Value = 1.4093975310840701
$ws.Range("H7").Style = $ws.Range("G7").Style
$ws.Range("I7").Value = 0.31039333343505798
